# Refresh crypto price/volume snapshot on Sheet1 (GitHub Actions scrape update).
# D column holds numeric-looking values stored as TEXT (e.g. "242.51"); we
# force the cells to Text format before writing so Excel keeps them as
# strings (matching the source data) instead of auto-converting to numbers.
# E column holds plain labels, written as normal strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$priceUpdates = @{
    "D2" = "242.51"
    "D3" = "21.52"
    "D4" = "5.242"
    "D5" = "0.05596"
    "D6" = "3.372"
    "D7" = "6.384"
    "D8" = "0.8070"
    "D9" = "0.9104"
    "D10" = "0.1426"
    "D11" = "0.07318"
    "D12" = "0.03208"
    "D13" = "0.03006"
    "D14" = "0.09275"
    "D15" = "3.609"
    "D16" = "0.001655"
    "D17" = "0.04703"
    "D18" = "0.0005802"
    "D19" = "0.006345"
    "D20" = "0.004981"
    "D21" = "0.001042"
    "D22" = "0.0001500"
    "D23" = "0.0003101"
    "D24" = "3.765"
    "D25" = "2.098"
    "D26" = "0.3273"
    "D40" = "0.03919"
    "D41" = "0.006930"
    "D42" = "0.003401"
    "D44" = "0.007484"
    "D45" = "0.00005945"
    "D46" = "0.00000000750"
    "D47" = "0.0005502"
    "D48" = "0.6826"
    "D49" = "0.06120"
    "D50" = "0.00002101"
    "D51" = "0.01010"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

$labelUpdates = @{
    "E18" = "17OneONE"
    "E47" = "46ACDXExchangeACXTWorstin24h"
}

foreach ($addr in $labelUpdates.Keys) {
    $ws.Range($addr).Value = $labelUpdates[$addr]
}
